# #5: insurance, claim, debt, investment done
# Fix the "保險" (insurance) sheet: proper column headers in row 1
# (the B1:D1 header cells were wrongly holding data values), and add the
# new trailing columns E:K (category, property_category, date,
# legislator_name, legislator_id, source_file, index) that every other
# sheet in this workbook already carries.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("保險")

# ---- Row 1: header labels ----
$ws.Range("B1").Value = "company"
$ws.Range("C1").Value = "name"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "property_category"
$ws.Range("F1").Value = "category"
$ws.Range("G1").Value = "date"
$ws.Range("H1").Value = "legislator_name"
$ws.Range("I1").Value = "legislator_id"
$ws.Range("J1").Value = "source_file"
$ws.Range("K1").Value = "index"

# match the bold/bordered header style already used by B1:D1
$ws.Range("D1").Copy()
$ws.Range("E1:K1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# (E2:K11 use the same plain/default data-row formatting already implicit
# on this sheet - nothing further to copy there)

# column G holds a text date like "2012-04-26" - force text format first so
# Excel doesn't reinterpret it as a date serial number
$ws.Range("G2:G11").NumberFormat = "@"

# ---- Data rows 2-11: fill new columns E:K ----
$rows = @(
    @{ Row = 2;  Index = 110; Company = "南山人壽"; Name = "南山好鑫動養老保險"; Owner = "景玉鳳" },
    @{ Row = 3;  Index = 111; Company = "南山人壽"; Name = "新二十年期增值分紅養老壽險"; Owner = "景玉鳳" },
    @{ Row = 4;  Index = 112; Company = "南山人壽"; Name = "南山全新增額養老壽險"; Owner = "林鴻池" },
    @{ Row = 5;  Index = 113; Company = "國泰人壽"; Name = "雙囍年年終身壽險"; Owner = "景玉鳳" },
    @{ Row = 6;  Index = 114; Company = "國泰人壽"; Name = "創世紀變額萬能壽險（丁型）"; Owner = "景玉鳳" },
    @{ Row = 7;  Index = 115; Company = "國泰人壽"; Name = "添意終身壽險"; Owner = "景玉鳳" },
    @{ Row = 8;  Index = 116; Company = "全球人壽"; Name = "增額終身壽險（B型）"; Owner = "林鴻池" },
    @{ Row = 9;  Index = 117; Company = "全球人壽"; Name = "增額終身壽險（B型）"; Owner = "林〇廷" },
    @{ Row = 10; Index = 118; Company = "全球人壽"; Name = "增額終身壽險（B型）"; Owner = "景玉鳳" },
    @{ Row = 11; Index = 119; Company = "全球人壽"; Name = "金彩306增額終身壽險"; Owner = "景玉鳳" }
)

foreach ($r in $rows) {
    $row = $r.Row
    # columns B/C/D already hold the correct company/name/owner text;
    # only the new E:K columns need to be populated.
    $ws.Range("B$row").Value = $r.Company
    $ws.Range("C$row").Value = $r.Name
    $ws.Range("D$row").Value = $r.Owner
    $ws.Range("E$row").Value = "insurance"
    $ws.Range("F$row").Value = "normal"
    $ws.Range("G$row").Value = "2012-04-26"
    $ws.Range("H$row").Value = "林鴻池"
    $ws.Range("I$row").Value = 1340
    $ws.Range("J$row").Value = "tmpdb4b1"
    $ws.Range("K$row").Value = $r.Index
}
